$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Range("B5").Value = 90
$ws.Range("B6").Value = 92
$ws.Range("B7").Value = 94
$ws.Range("B8").Value = 95
$ws.Range("B9").Value = 97
$ws.Range("B10").Value = 99
$ws.Range("B11").Value = 101
$ws.Range("B12").Value = 103
$ws.Range("B13").Value = 105
$ws.Range("B14").Value = 106
$ws.Range("B15").Value = 108
$ws.Range("B16").Value = 110
$ws.Range("B17").Value = 112
$ws.Range("B18").Value = 114
$ws.Range("B19").Value = 116
$ws.Range("B20").Value = 117
$ws.Range("B21").Value = 119
$ws.Range("B22").Value = 121
$ws.Range("B23").Value = 123
$ws.Range("B24").Value = 126
$ws.Range("B25").Value = 128

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Range("B3").Value = 77
$ws.Range("B4").Value = 79
$ws.Range("B5").Value = 81
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 85
$ws.Range("B8").Value = 87
$ws.Range("B9").Value = 89
$ws.Range("B10").Value = 91
$ws.Range("B11").Value = 93
$ws.Range("B12").Value = 95
$ws.Range("B13").Value = 97
$ws.Range("B16").Value = 102
$ws.Range("B17").Value = 104
$ws.Range("B22").Value = 112
$ws.Range("B23").Value = 114
$ws.Range("B24").Value = 116
$ws.Range("B25").Value = 118
$ws.Range("B26").Value = 120
$ws.Range("B27").Value = 122
$ws.Range("B28").Value = 124
$ws.Range("B29").Value = 126
$ws.Range("B30").Value = 129

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Range("B2").Value = 66
$ws.Range("B3").Value = 68
$ws.Range("B7").Value = 78
$ws.Range("B8").Value = 80
$ws.Range("B9").Value = 82
$ws.Range("B10").Value = 84
$ws.Range("B11").Value = 86
$ws.Range("B12").Value = 88
$ws.Range("B13").Value = 90
$ws.Range("B14").Value = 91
$ws.Range("B15").Value = 93
$ws.Range("B16").Value = 95
$ws.Range("B17").Value = 97
$ws.Range("B18").Value = 99
$ws.Range("B19").Value = 100
$ws.Range("B20").Value = 102
$ws.Range("B21").Value = 104
$ws.Range("B22").Value = 106
$ws.Range("B31").Value = 124

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 60
$ws.Range("B4").Value = 63
$ws.Range("B5").Value = 65
$ws.Range("B6").Value = 68
$ws.Range("B7").Value = 70
$ws.Range("B8").Value = 72
$ws.Range("B9").Value = 75
$ws.Range("B10").Value = 77
$ws.Range("B11").Value = 79
$ws.Range("B12").Value = 81
$ws.Range("B13").Value = 83
$ws.Range("B14").Value = 85
$ws.Range("B15").Value = 87
$ws.Range("B16").Value = 89
$ws.Range("B17").Value = 90
$ws.Range("B18").Value = 92
$ws.Range("B19").Value = 94
$ws.Range("B20").Value = 96
$ws.Range("B21").Value = 98
$ws.Range("B22").Value = 100
$ws.Range("B23").Value = 102
$ws.Range("B24").Value = 103
$ws.Range("B25").Value = 105
$ws.Range("B26").Value = 107
$ws.Range("B27").Value = 109
$ws.Range("B28").Value = 112
$ws.Range("B29").Value = 114
$ws.Range("B30").Value = 116
$ws.Range("B31").Value = 119
$ws.Range("B32").Value = 122

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Range("B2").Value = 50
$ws.Range("B3").Value = 53
$ws.Range("B4").Value = 56
$ws.Range("B5").Value = 58
$ws.Range("B6").Value = 61
$ws.Range("B7").Value = 63
$ws.Range("B8").Value = 66
$ws.Range("B9").Value = 68
$ws.Range("B10").Value = 70
$ws.Range("B11").Value = 72
$ws.Range("B25").Value = 100
$ws.Range("B26").Value = 102
$ws.Range("B28").Value = 107
$ws.Range("B29").Value = 110
$ws.Range("B30").Value = 112
$ws.Range("B32").Value = 121

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Range("B2").Value = 43
$ws.Range("B3").Value = 46
$ws.Range("B4").Value = 49
$ws.Range("B5").Value = 52
$ws.Range("B6").Value = 54
$ws.Range("B7").Value = 57
$ws.Range("B8").Value = 59
$ws.Range("B9").Value = 61
$ws.Range("B10").Value = 63
$ws.Range("B11").Value = 66
$ws.Range("B12").Value = 68
$ws.Range("B13").Value = 70
$ws.Range("B14").Value = 72
$ws.Range("B15").Value = 74
$ws.Range("B16").Value = 76
$ws.Range("B17").Value = 78
$ws.Range("B18").Value = 80
$ws.Range("B19").Value = 82
$ws.Range("B20").Value = 84
$ws.Range("B21").Value = 86
$ws.Range("B22").Value = 88
$ws.Range("B23").Value = 91
$ws.Range("B24").Value = 93
$ws.Range("B25").Value = 95
$ws.Range("B26").Value = 97
$ws.Range("B27").Value = 100
$ws.Range("B28").Value = 103
$ws.Range("B29").Value = 106
$ws.Range("B30").Value = 109
$ws.Range("B31").Value = 117
$ws.Range("B32").Value = 117
